$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = '[ You%Li%null%1,              Rachel M%Reeves%null%1,              Xin%Wang%null%1,              Quique%Bassat%null%1,              W Abdullah%Brooks%null%1,              Cheryl%Cohen%null%1,              David P%Moore%null%1,              Marta%Nunes%null%1,              Barbara%Rath%null%1,              Harry%Campbell%null%1,              Harish%Nair%null%1,              Sozinho%Acacio%null%1,              Wladimir J%Alonso%null%1,              Martin%Antonio%null%1,              Guadalupe%Ayora Talavera%null%1,              Darmaa%Badarch%null%1,              Vicky L%Baillie%null%1,              Gisela%Barrera-Badillo%null%1,              Godfrey%Bigogo%null%1,              Shobha%Broor%null%1,              Dana%Bruden%null%1,              Philippe%Buchy%null%1,              Peter%Byass%null%1,              James%Chipeta%null%1,              Wilfrido%Clara%null%1,              Duc-Anh%Dang%null%1,              Carla Cecília%de Freitas Lázaro Emediato%null%1,              Menno%de Jong%null%1,              José Alberto%Díaz-Quiñonez%null%1,              Lien Anh Ha%Do%null%1,              Rodrigo A%Fasce%null%1,              Luzhao%Feng%null%1,              Mark J%Ferson%null%1,              Angela%Gentile%null%1,              Bradford D%Gessner%null%1,              Doli%Goswami%null%1,              Sophie%Goyet%null%1,              Carlos G%Grijalva%null%1,              Natasha%Halasa%null%1,              Orienka%Hellferscee%null%1,              Danielle%Hessong%null%1,              Nusrat%Homaira%null%1,              Jorge%Jara%null%1,              Kathleen%Kahn%null%1,              Najwa%Khuri-Bulos%null%1,              Karen L%Kotloff%null%1,              Claudio F%Lanata%null%1,              Olga%Lopez%null%1,              Maria Renee%Lopez Bolaños%null%1,              Marilla G%Lucero%null%1,              Florencia%Lucion%null%1,              Socorro P%Lupisan%null%1,              Shabir A%Madhi%null%1,              Omphile%Mekgoe%null%1,              Cinta%Moraleda%null%1,              Jocelyn%Moyes%null%1,              Kim%Mulholland%null%1,              Patrick K%Munywoki%null%1,              Fathima%Naby%null%1,              Thanh Hung%Nguyen%null%1,              Mark P%Nicol%null%1,              D James%Nokes%null%1,              Daniel E%Noyola%null%1,              Daisuke%Onozuka%null%1,              Nandhini%Palani%null%1,              Yong%Poovorawan%null%1,              Mustafizur%Rahman%null%1,              Kaat%Ramaekers%null%1,              Candice%Romero%null%1,              Elizabeth P%Schlaudecker%null%1,              Brunhilde%Schweiger%null%1,              Phil%Seidenberg%null%1,              Eric A F%Simoes%null%1,              Rosalyn%Singleton%null%1,              Sujatha%Sistla%null%1,              Katharine%Sturm-Ramirez%null%1,              Nungruthai%Suntronwong%null%1,              Agustinus%Sutanto%null%1,              Milagritos D%Tapia%null%1,              Somsak%Thamthitiwat%null%1,              Ilada%Thongpan%null%1,              Gayani%Tillekeratne%null%1,              Yeny O%Tinoco%null%1,              Florette K%Treurnicht%null%1,              Claudia%Turner%null%1,              Paul%Turner%null%1,              Rogier%van Doorn%null%1,              Marc%Van Ranst%null%2,              Benoit%Visseaux%null%1,              Sunthareeya%Waicharoen%null%1,              Jianwei%Wang%null%1,              Lay-Myint%Yoshida%null%1,              Heather J%Zar%null%1,          You%Li%null%1,          Rachel M%Reeves%null%1,          Xin%Wang%null%1,          Quique%Bassat%null%1,          W Abdullah%Brooks%null%1,          Cheryl%Cohen%null%1,          David P%Moore%null%1,          Marta%Nunes%null%1,          Barbara%Rath%null%1,          Harry%Campbell%null%1,          Harish%Nair%null%1,          Sozinho%Acacio%null%1,          Wladimir J%Alonso%null%1,          Martin%Antonio%null%1,          Guadalupe%Ayora Talavera%null%1,          Darmaa%Badarch%null%1,          Vicky L%Baillie%null%1,          Gisela%Barrera-Badillo%null%1,          Godfrey%Bigogo%null%1,          Shobha%Broor%null%1,          Dana%Bruden%null%1,          Philippe%Buchy%null%1,          Peter%Byass%null%1,          James%Chipeta%null%1,          Wilfrido%Clara%null%1,          Duc-Anh%Dang%null%1,          Carla Cecília%de Freitas Lázaro Emediato%null%1,          Menno%de Jong%null%1,          José Alberto%Díaz-Quiñonez%null%1,          Lien Anh Ha%Do%null%1,          Rodrigo A%Fasce%null%1,          Luzhao%Feng%null%2,          Mark J%Ferson%null%1,          Angela%Gentile%null%1,          Bradford D%Gessner%null%1,          Doli%Goswami%null%1,          Sophie%Goyet%null%1,          Carlos G%Grijalva%null%1,          Natasha%Halasa%null%1,          Orienka%Hellferscee%null%1,          Danielle%Hessong%null%1,          Nusrat%Homaira%null%1,          Jorge%Jara%null%1,          Kathleen%Kahn%null%1,          Najwa%Khuri-Bulos%null%1,          Karen L%Kotloff%null%1,          Claudio F%Lanata%null%1,          Olga%Lopez%null%1,          Maria Renee%Lopez Bolaños%null%1,          Marilla G%Lucero%null%1,          Florencia%Lucion%null%1,          Socorro P%Lupisan%null%1,          Shabir A%Madhi%null%1,          Omphile%Mekgoe%null%1,          Cinta%Moraleda%null%1,          Jocelyn%Moyes%null%1,          Kim%Mulholland%null%1,          Patrick K%Munywoki%null%1,          Fathima%Naby%null%1,          Thanh Hung%Nguyen%null%1,          Mark P%Nicol%null%1,          D James%Nokes%null%1,          Daniel E%Noyola%null%1,          Daisuke%Onozuka%null%1,          Nandhini%Palani%null%1,          Yong%Poovorawan%null%1,          Mustafizur%Rahman%null%1,          Kaat%Ramaekers%null%1,          Candice%Romero%null%1,          Elizabeth P%Schlaudecker%null%1,          Brunhilde%Schweiger%null%1,          Phil%Seidenberg%null%1,          Eric A F%Simoes%null%1,          Rosalyn%Singleton%null%1,          Sujatha%Sistla%null%1,          Katharine%Sturm-Ramirez%null%1,          Nungruthai%Suntronwong%null%1,          Agustinus%Sutanto%null%1,          Milagritos D%Tapia%null%1,          Somsak%Thamthitiwat%null%1,          Ilada%Thongpan%null%1,          Gayani%Tillekeratne%null%1,          Yeny O%Tinoco%null%1,          Florette K%Treurnicht%null%1,          Claudia%Turner%null%1,          Paul%Turner%null%1,          Rogier%van Doorn%null%1,          Marc%Van Ranst%null%3,          Benoit%Visseaux%null%1,          Sunthareeya%Waicharoen%null%1,          Jianwei%Wang%null%5,          Lay-Myint%Yoshida%null%1,          Heather J%Zar%null%1]'
$ws.Range("E3").Value2 = '[Adeeba%Kamarulzaman%coreGivesNoEmail%1,            Jack%Bee Chook%coreGivesNoEmail%1,            Kim%Tien Ng%coreGivesNoEmail%1,            Kok%Keng Tee%coreGivesNoEmail%1,            Maryam%Nabiel Al-Khannaq%coreGivesNoEmail%1,            Nik%Sherina Hanafi%coreGivesNoEmail%1,            van%der Hoek%coreGivesNoEmail%2,            van%der Hoek%coreGivesNoEmail%0,            van%Elden%coreGivesNoEmail%1,            Xiang%Yong Oong%coreGivesNoEmail%1,            Yong%Kek Pang%coreGivesNoEmail%1,            Yutaka%Takebe%coreGivesNoEmail%1]'
$ws.Range("E4").Value2 = '[Robin%Brittain-Long%xref no email%1,    Lars-Magnus%Andersson%xref no email%1,    Sigvard%Olofsson%xref no email%1,    Magnus%Lindh%xref no email%1,    Johan%Westin%xref no email%1]'
$ws.Range("E5").Value2 = '[Binglin%Cui%binglincui@stu.edu.cn%1,                Dangui%Zhang%danguizhang@stu.edu.cn%1,                Hui%Pan%michaelpan@stu.edu.cn%1,                Fan%Zhang%lemon_fan@163.com%1,                Jeremy%Farrar%jeremyjfarrar@gmail.com%1,                Frieda%Law%lawf@stu.edu.cn%1,                H Rogier%van Doorn%hrogier@gmail.com%1,                Beiyan%Wu%1261305798@qq.com%1,                William%Ba-Thein%wbathein@stu.edu.cn%1,               Binglin%Cui%binglincui@stu.edu.cn%1,               Dangui%Zhang%danguizhang@stu.edu.cn%1,               Hui%Pan%michaelpan@stu.edu.cn%1,               Fan%Zhang%lemon_fan@163.com%1,               Jeremy%Farrar%jeremyjfarrar@gmail.com%1,               Frieda%Law%lawf@stu.edu.cn%1,               H Rogier%van Doorn%hrogier@gmail.com%1,               Beiyan%Wu%1261305798@qq.com%1,               William%Ba-Thein%wbathein@stu.edu.cn%1]'
$ws.Range("E6").Value2 = '[Amy%Ransier%coreGivesNoEmail%1,            Andrew%Berglund%coreGivesNoEmail%1,            Christine%C. Robinson%coreGivesNoEmail%1,            Christopher%D. Town%coreGivesNoEmail%1,            Edison%Luiz Durigon%coreGivesNoEmail%1,            Jos\u00e9%Antonio Jerez%coreGivesNoEmail%1,            Kathryn%V. Holmes%coreGivesNoEmail%1,            Luiz%Gustavo Bentim G\u00f3es%coreGivesNoEmail%1,            Nadia%Fedorova%coreGivesNoEmail%1,            Philip%A. Weston%coreGivesNoEmail%1,            Rebecca%A. Halpin%coreGivesNoEmail%1,            Samuel%R. Dominguez%coreGivesNoEmail%1,            Susmita%Shrivastava%coreGivesNoEmail%1,            Zhaohui%Qian%coreGivesNoEmail%1]'
$ws.Range("E7").Value2 = '[Yifei%Fu%NULL%1,                Lifeng%Pan%NULL%1,                Qiao%Sun%NULL%1,                Weiping%Zhu%NULL%1,                Linying%Zhu%NULL%1,                Chuchu%Ye%NULL%1,                Caoyi%Xue%NULL%1,                Yuanping%Wang%NULL%1,                Qing%Liu%NULL%1,                Ping%Ma%NULL%1,                Huifang%Qiu%NULL%1,                Srinand%Sreevatsan%NULL%1,               Yifei%Fu%NULL%1,               Lifeng%Pan%NULL%1,               Qiao%Sun%NULL%1,               Weiping%Zhu%NULL%1,               Linying%Zhu%NULL%1,               Chuchu%Ye%NULL%1,               Caoyi%Xue%NULL%1,               Yuanping%Wang%NULL%1,               Qing%Liu%NULL%1,               Ping%Ma%NULL%1,               Huifang%Qiu%NULL%1,               Srinand%Sreevatsan%NULL%1]'
$ws.Range("E8").Value2 = '[Safak%Goktas%NULL%1,                Mumtaz Cem%Sirin%NULL%1,               Safak%Goktas%NULL%1,               Mumtaz Cem%Sirin%NULL%1]'
$ws.Range("E9").Value2 = '[Inger%Heimdal%NULL%1,                Nina%Moe%NULL%1,                Sidsel%Krokstad%NULL%2,                Andreas%Christensen%NULL%2,                Lars Høsøien%Skanke%NULL%2,                Svein Arne%Nordbø%NULL%2,                Henrik%Døllner%henrik.dollner@ntnu.no%2,               Inger%Heimdal%NULL%1,               Nina%Moe%NULL%2,               Sidsel%Krokstad%NULL%2,               Andreas%Christensen%NULL%2,               Lars Høsøien%Skanke%NULL%2,               Svein Arne%Nordbø%NULL%2,               Henrik%Døllner%henrik.dollner@ntnu.no%2]'
$ws.Range("E10").Value2 = '[Nina%Moe%NULL%1,                Inger Heimdal%Stenseng%NULL%1,                Sidsel%Krokstad%NULL%0,                Andreas%Christensen%NULL%0,                Lars Høsøien%Skanke%NULL%0,                Kari Ravndal%Risnes%NULL%1,                Svein Arne%Nordbø%NULL%0,                Henrik%Døllner%henrik.dollner@ntnu.no%0,               Nina%Moe%NULL%0,               Inger Heimdal%Stenseng%NULL%1,               Sidsel%Krokstad%NULL%0,               Andreas%Christensen%NULL%0,               Lars Høsøien%Skanke%NULL%0,               Kari Ravndal%Risnes%NULL%1,               Svein Arne%Nordbø%NULL%0,               Henrik%Døllner%henrik.dollner@ntnu.no%0]'
$ws.Range("E11").Value2 = '[Seema%Jain%xref no email%1,    Derek J.%Williams%xref no email%1,    Sandra R.%Arnold%xref no email%1,    Krow%Ampofo%xref no email%1,    Anna M.%Bramley%xref no email%1,    Carrie%Reed%xref no email%1,    Chris%Stockmann%xref no email%1,    Evan J.%Anderson%xref no email%1,    Carlos G.%Grijalva%xref no email%1,    Wesley H.%Self%xref no email%1,    Yuwei%Zhu%xref no email%1,    Anami%Patel%xref no email%1,    Weston%Hymas%xref no email%1,    James D.%Chappell%xref no email%1,    Robert A.%Kaufman%xref no email%1,    J. Herman%Kan%xref no email%1,    David%Dansie%xref no email%1,    Noel%Lenny%xref no email%1,    David R.%Hillyard%xref no email%1,    Lia M.%Haynes%xref no email%1,    Min%Levine%xref no email%1,    Stephen%Lindstrom%xref no email%1,    Jonas M.%Winchell%xref no email%1,    Jacqueline M.%Katz%xref no email%1,    Dean%Erdman%xref no email%1,    Eileen%Schneider%xref no email%1,    Lauri A.%Hicks%xref no email%1,    Richard G.%Wunderink%xref no email%1,    Kathryn M.%Edwards%xref no email%1,    Andrew T.%Pavia%xref no email%1,    Jonathan A.%McCullers%xref no email%1,    Lyn%Finelli%xref no email%1]'
$ws.Range("E12").Value2 = '[Monika%Jevšnik%NULL%1,                Andrej%Steyer%NULL%1,                Marko%Pokorn%NULL%1,                Tatjana%Mrvič%NULL%1,                Štefan%Grosek%NULL%1,                Franc%Strle%NULL%1,                Lara%Lusa%NULL%1,                Miroslav%Petrovec%NULL%1,                Oliver%Schildgen%NULL%1,               Monika%Jevšnik%NULL%1,               Andrej%Steyer%NULL%1,               Marko%Pokorn%NULL%1,               Tatjana%Mrvič%NULL%1,               Štefan%Grosek%NULL%1,               Franc%Strle%NULL%1,               Lara%Lusa%NULL%1,               Miroslav%Petrovec%NULL%1,               Oliver%Schildgen%NULL%1]'
$ws.Range("E13").Value2 = '[Eeva K%Broberg%NULL%1,                Matti%Waris%NULL%1,                Kari%Johansen%NULL%1,                René%Snacken%NULL%1,                Pasi%Penttinen%NULL%1,                NULL%NULL%NULL%3,               Eeva K%Broberg%NULL%1,               Matti%Waris%NULL%1,               Kari%Johansen%NULL%1,               René%Snacken%NULL%1,               Pasi%Penttinen%NULL%1,               NULL%NULL%NULL%2]'
$ws.Range("E14").Value2 = '[Xiongfei%Ju%NULL%1,                Qiaoyun%Fang%NULL%1,                Jian%Zhang%NULL%1,                Angao%Xu%NULL%1,                Lihuan%Liang%NULL%1,                Changwen%Ke%hzmedicine@163.com%1,               Xiongfei%Ju%NULL%1,               Qiaoyun%Fang%NULL%1,               Jian%Zhang%NULL%1,               Angao%Xu%NULL%1,               Lihuan%Liang%NULL%1,               Changwen%Ke%hzmedicine@163.com%1]'
$ws.Range("E15").Value2 = '[Marie E.%Killerby%lxo9@cdc.gov%1,                Holly M.%Biggs%NULL%1,                Amber%Haynes%NULL%1,                Rebecca M.%Dahl%NULL%1,                Desiree%Mustaquim%NULL%1,                Susan I.%Gerber%NULL%1,                John T.%Watson%NULL%1,               Marie E.%Killerby%lxo9@cdc.gov%1,               Holly M.%Biggs%NULL%1,               Amber%Haynes%NULL%1,               Rebecca M.%Dahl%NULL%1,               Desiree%Mustaquim%NULL%1,               Susan I.%Gerber%NULL%1,               John T.%Watson%NULL%1]'
$ws.Range("E16").Value2 = '[Jeong‐Min%Kim%NULL%1,                Hee‐Dong%Jung%NULL%1,                Hyang‐Min%Cheong%NULL%1,                Anna%Lee%NULL%1,                Nam‐Joo%Lee%NULL%1,                Hyuk%Chu%NULL%1,                Joo‐Yeon%Lee%NULL%1,                Sung Soon%Kim%NULL%1,                Jang‐Hoon%Choi%jhchoi@nih.go.kr%1,               Jeong‐Min%Kim%NULL%1,               Hee‐Dong%Jung%NULL%1,               Hyang‐Min%Cheong%NULL%1,               Anna%Lee%NULL%1,               Nam‐Joo%Lee%NULL%1,               Hyuk%Chu%NULL%1,               Joo‐Yeon%Lee%NULL%1,               Sung Soon%Kim%NULL%1,               Jang‐Hoon%Choi%jhchoi@nih.go.kr%1]'
$ws.Range("E17").Value2 = '[Susanna K. P.%Lau%xref no email%1,    Patrick C. Y.%Woo%xref no email%1,    Cyril C. Y.%Yip%xref no email%1,    Herman%Tse%xref no email%1,    Hoi-wah%Tsoi%xref no email%1,    Vincent C. C.%Cheng%xref no email%1,    Paul%Lee%xref no email%1,    Bone S. F.%Tang%xref no email%1,    Chris H. Y.%Cheung%xref no email%1,    Rodney A.%Lee%xref no email%1,    Lok-yee%So%xref no email%1,    Yu-lung%Lau%xref no email%1,    Kwok-hung%Chan%xref no email%1,    Kwok-yung%Yuen%xref no email%1]'
$ws.Range("E18").Value2 = '[Prawin%Kumar%xref no email%1,    Guruprasad R.%Medigeshi%xref no email%1,    Vishnu S.%Mishra%xref no email%1,    Mojahidul%Islam%xref no email%1,    Shivani%Randev%xref no email%1,    Aparna%Mukherjee%xref no email%1,    Rama%Chaudhry%xref no email%1,    Arti%Kapil%xref no email%1,    Kana%Ram Jat%xref no email%1,    Rakesh%Lodha%xref no email%1,    Sushil K.%Kabra%xref no email%1]'
$ws.Range("E19").Value2 = '[Seung Jin%Lee%xref no email%1,    Shin Hae%Lee%xref no email%1,    Eun Kyo%Ha%xref no email%1,    Youn Ho%Sheen%xref no email%1,    Myong Soon%Sung%xref no email%1,    Young-Ho%Jung%xref no email%1,    Kyung Suk%Lee%xref no email%1,    Hye Mi%Jee%xref no email%1,    Man Yong%Han%xref no email%1]'
$ws.Range("E20").Value2 = '[Xuechao%Li%NULL%1,                Juansheng%Li%lijsh16@163.com%1,                Lei%Meng%mleicdc@163.com%1,                Wanqi%Zhu%NULL%1,                Xinfeng%Liu%NULL%1,                Mei%Yang%NULL%1,                Deshan%Yu%NULL%2,                Lixia%Niu%NULL%1,                Xiping%Shen%NULL%1,               Xuechao%Li%NULL%1,               Juansheng%Li%lijsh16@163.com%1,               Lei%Meng%mleicdc@163.com%1,               Wanqi%Zhu%NULL%1,               Xinfeng%Liu%NULL%1,               Mei%Yang%NULL%1,               Deshan%Yu%NULL%2,               Lixia%Niu%NULL%1,               Xiping%Shen%NULL%1]'
$ws.Range("E21").Value2 = '[Christine M.%Litwin%clitwin@gru.edu%1,                James G.%Bosley%NULL%1,               Christine M.%Litwin%clitwin@gru.edu%1,               James G.%Bosley%NULL%1]'
$ws.Range("E22").Value2 = '[Yohei%Matoba%xref no email%1,    Chieko%Abiko%xref no email%1,    Tatsuya%Ikeda%xref no email%1,    Yoko%Aoki%xref no email%1,    Yu%Suzuki%xref no email%1,    Kazue%Yahagi%xref no email%1,    Yoko%Matsuzaki%xref no email%1,    Tsutomu%Itagaki%xref no email%1,    Fumio%Katsushima%xref no email%1,    Yuriko%Katsushima%xref no email%1,    Katsumi%Mizuta%xref no email%1]'
$ws.Range("E23").Value2 = '[Richard%Njouom%NULL%1,                Elsie Laban%Yekwa%NULL%1,                Pierre%Cappy%NULL%1,                Astrid%Vabret%NULL%2,                Pascal%Boisier%NULL%1,                Dominique%Rousset%NULL%1,               Richard%Njouom%NULL%1,               Elsie Laban%Yekwa%NULL%1,               Pierre%Cappy%NULL%1,               Astrid%Vabret%NULL%2,               Pascal%Boisier%NULL%1,               Dominique%Rousset%NULL%1]'
$ws.Range("E24").Value2 = '[Minna%Paloniemi%NULL%1,                Suvi%Lappalainen%NULL%1,                Timo%Vesikari%NULL%1,               Minna%Paloniemi%NULL%1,               Suvi%Lappalainen%NULL%1,               Timo%Vesikari%NULL%1]'
$ws.Range("E25").Value2 = '[Lili%Ren%NULL%1,                Richard%Gonzalez%NULL%1,                Jin%Xu%NULL%1,                Yan%Xiao%NULL%1,                Yongjun%Li%NULL%1,                Hongli%Zhou%NULL%1,                Jianguo%Li%NULL%1,                Qingqing%Yang%NULL%1,                Jing%Zhang%NULL%1,                Lan%Chen%NULL%1,                Wei%Wang%NULL%1,                Guy%Vernet%NULL%1,                Gláucia%Paranhos‐Baccalà%NULL%1,                Zhong%Wang%wangzhong523@vip.163.com%1,                Jianwei%Wang%wangjw28@ipbcams.ac.cn%1,               Lili%Ren%NULL%1,               Richard%Gonzalez%NULL%1,               Jin%Xu%NULL%1,               Yan%Xiao%NULL%1,               Yongjun%Li%NULL%1,               Hongli%Zhou%NULL%1,               Jianguo%Li%NULL%1,               Qingqing%Yang%NULL%1,               Jing%Zhang%NULL%1,               Lan%Chen%NULL%1,               Wei%Wang%NULL%1,               Guy%Vernet%NULL%1,               Gláucia%Paranhos‐Baccalà%NULL%1,               Zhong%Wang%wangzhong523@vip.163.com%1,               Jianwei%Wang%wangjw28@ipbcams.ac.cn%1]'
$ws.Range("E26").Value2 = '[Albert%Jan%coreGivesNoEmail%1,            Druelle%Valentin%coreGivesNoEmail%1,            Dyrdak%Robert%coreGivesNoEmail%1,            Hodcroft%Emma B.%coreGivesNoEmail%1,            Neher%Richard A.%coreGivesNoEmail%1]'
$ws.Range("E27").Value2 = '[Rapeepun%Soonnarong%sakura.cu@hotmail.com%1,                Ilada%Thongpan%ilada.cu09@gmail.com%1,                Sunchai%Payungporn%sp.medbiochemcu@gmail.com%1,                Chanpim%Vuthitanachot%viboovut@hotmail.com%1,                Viboonsuk%Vuthitanachot%viboovut@hotmail.com%1,                Preeyaporn%Vichiwattana%preeya_teiy@hotmail.com%1,                Sompong%Vongpunsawad%sompong.vo@chula.ac.th%1,                Yong%Poovorawan%yong.p@chula.ac.th%1,               Rapeepun%Soonnarong%sakura.cu@hotmail.com%1,               Ilada%Thongpan%ilada.cu09@gmail.com%1,               Sunchai%Payungporn%sp.medbiochemcu@gmail.com%1,               Chanpim%Vuthitanachot%viboovut@hotmail.com%1,               Viboonsuk%Vuthitanachot%viboovut@hotmail.com%1,               Preeyaporn%Vichiwattana%preeya_teiy@hotmail.com%1,               Sompong%Vongpunsawad%sompong.vo@chula.ac.th%1,               Yong%Poovorawan%yong.p@chula.ac.th%1]'
$ws.Range("E29").Value2 = '[S M Iftekhar%Uddin%NULL%1,                Janet A%Englund%NULL%1,                Jane Y%Kuypers%NULL%1,                Helen Y%Chu%NULL%1,                Mark C%Steinhoff%NULL%1,                Subarna K%Khatry%NULL%1,                Steve C%LeClerq%NULL%1,                James M%Tielsch%NULL%1,                Luke C%Mullany%NULL%1,                Laxman%Shrestha%NULL%2,                Joanne%Katz%jkatz1@jhu.edu%2,               S M Iftekhar%Uddin%NULL%1,               Janet A%Englund%NULL%1,               Jane Y%Kuypers%NULL%1,               Helen Y%Chu%NULL%1,               Mark C%Steinhoff%NULL%1,               Subarna K%Khatry%NULL%1,               Steve C%LeClerq%NULL%1,               James M%Tielsch%NULL%1,               Luke C%Mullany%NULL%1,               Laxman%Shrestha%NULL%2,               Joanne%Katz%jkatz1@jhu.edu%2]'
$ws.Range("E30").Value2 = '[Helen Y.%Chu%helenchu@uw.edu%1,                Joanne%Katz%NULL%0,                James%Tielsch%NULL%1,                Subarna K.%Khatry%NULL%1,                Laxman%Shrestha%NULL%0,                Steven C.%LeClerq%NULL%1,                Amalia%Magaret%NULL%1,                Jane%Kuypers%NULL%1,                Mark%Steinhoff%NULL%1,                Janet A.%Englund%NULL%1,               Helen Y.%Chu%helenchu@uw.edu%1,               Joanne%Katz%NULL%0,               James%Tielsch%NULL%1,               Subarna K.%Khatry%NULL%1,               Laxman%Shrestha%NULL%0,               Steven C.%LeClerq%NULL%1,               Amalia%Magaret%NULL%1,               Jane%Kuypers%NULL%1,               Mark%Steinhoff%NULL%1,               Janet A.%Englund%NULL%1]'
$ws.Range("E31").Value2 = '[Lia%van der Hoek%NULL%1,                Gabriele%Ihorst%NULL%1,                Klaus%Sure%NULL%1,                Astrid%Vabret%NULL%0,                Ronald%Dijkman%NULL%1,                Michel%de Vries%NULL%1,                Johannes%Forster%NULL%1,                Ben%Berkhout%NULL%1,                Klaus%Uberla%NULL%1,               Lia%van der Hoek%NULL%1,               Gabriele%Ihorst%NULL%1,               Klaus%Sure%NULL%1,               Astrid%Vabret%NULL%0,               Ronald%Dijkman%NULL%1,               Michel%de Vries%NULL%1,               Johannes%Forster%NULL%1,               Ben%Berkhout%NULL%1,               Klaus%Uberla%NULL%1]'
$ws.Range("E32").Value2 = '[Litty%Varghese%NULL%1,                Philip%Zachariah%pz2177@cumc.columbia.edu%1,                Celibell%Vargas%NULL%1,                Philip%LaRussa%NULL%1,                Ryan T%Demmer%NULL%1,                Yoko E%Furuya%NULL%1,                Susan%Whittier%NULL%1,                Carrie%Reed%NULL%1,                Melissa S%Stockwell%NULL%1,                Lisa%Saiman%NULL%1,               Litty%Varghese%NULL%1,               Philip%Zachariah%pz2177@cumc.columbia.edu%1,               Celibell%Vargas%NULL%1,               Philip%LaRussa%NULL%1,               Ryan T%Demmer%NULL%1,               Yoko E%Furuya%NULL%1,               Susan%Whittier%NULL%1,               Carrie%Reed%NULL%1,               Melissa S%Stockwell%NULL%1,               Lisa%Saiman%NULL%1]'
$ws.Range("E33").Value2 = '[Marietjie%Venter%marietjiev@nicd.ac.za%1,                Ria%Lassaunière%NULL%1,                Tina Louise%Kresfelder%NULL%1,                Yvette%Westerberg%NULL%1,                Adele%Visser%NULL%1,               Marietjie%Venter%marietjiev@nicd.ac.za%1,               Ria%Lassaunière%NULL%1,               Tina Louise%Kresfelder%NULL%1,               Yvette%Westerberg%NULL%1,               Adele%Visser%NULL%1]'
$ws.Range("E34").Value2 = '[Benoit%Visseaux%NULL%1,                Charles%Burdet%NULL%1,                Guillaume%Voiriot%NULL%1,                François-Xavier%Lescure%NULL%1,                Taous%Chougar%NULL%1,                Olivier%Brugière%NULL%1,                Bruno%Crestani%NULL%1,                Enrique%Casalino%NULL%1,                Charlotte%Charpentier%NULL%1,                Diane%Descamps%NULL%1,                Jean-François%Timsit%NULL%1,                Yazdan%Yazdanpanah%NULL%1,                Nadhira%Houhou-Fidouh%NULL%1,                Dena L.%Schanzer%NULL%1,               Benoit%Visseaux%NULL%1,               Charles%Burdet%NULL%1,               Guillaume%Voiriot%NULL%1,               François-Xavier%Lescure%NULL%1,               Taous%Chougar%NULL%1,               Olivier%Brugière%NULL%1,               Bruno%Crestani%NULL%1,               Enrique%Casalino%NULL%1,               Charlotte%Charpentier%NULL%1,               Diane%Descamps%NULL%1,               Jean-François%Timsit%NULL%1,               Yazdan%Yazdanpanah%NULL%1,               Nadhira%Houhou-Fidouh%NULL%1,               Dena L.%Schanzer%NULL%1]'
$ws.Range("E35").Value2 = '[Guohong%Huang%NULL%1,                Deshan%Yu%NULL%0,                Naiying%Mao%NULL%1,                Zhen%Zhu%NULL%1,                Hui%Zhang%NULL%1,                Zhongyi%Jiang%NULL%1,                Hongyu%Li%NULL%1,                Yan%Zhang%NULL%1,                Jing%Shi%NULL%1,                Shuang%Zhang%NULL%1,                Xinhua%Wang%NULL%1,                Wenbo%Xu%NULL%1,                Dong-Yan%Jin%NULL%1,               Guohong%Huang%NULL%1,               Deshan%Yu%NULL%0,               Naiying%Mao%NULL%1,               Zhen%Zhu%NULL%1,               Hui%Zhang%NULL%1,               Zhongyi%Jiang%NULL%1,               Hongyu%Li%NULL%1,               Yan%Zhang%NULL%1,               Jing%Shi%NULL%1,               Shuang%Zhang%NULL%1,               Xinhua%Wang%NULL%1,               Wenbo%Xu%NULL%1,               Dong-Yan%Jin%NULL%1]'
$ws.Range("E36").Value2 = '[Su-fen%Zhang%NULL%1,                Jiu-ling%Tuo%NULL%1,                Xu-bin%Huang%NULL%1,                Xun%Zhu%NULL%1,                Ding-mei%Zhang%NULL%1,                Kai%Zhou%NULL%1,                Lei%Yuan%NULL%1,                Hong-jiao%Luo%NULL%1,                Bo-jian%Zheng%NULL%1,                Kwok-yung%Yuen%NULL%1,                Meng-feng%Li%NULL%1,                Kai-yuan%Cao%NULL%1,                Lin%Xu%NULL%1,                Stefan%Pöhlmann%NULL%1,               Su-fen%Zhang%NULL%1,               Jiu-ling%Tuo%NULL%1,               Xu-bin%Huang%NULL%1,               Xun%Zhu%NULL%1,               Ding-mei%Zhang%NULL%1,               Kai%Zhou%NULL%1,               Lei%Yuan%NULL%1,               Hong-jiao%Luo%NULL%1,               Bo-jian%Zheng%NULL%1,               Kwok-yung%Yuen%NULL%1,               Meng-feng%Li%NULL%1,               Kai-yuan%Cao%NULL%1,               Lin%Xu%NULL%1,               Stefan%Pöhlmann%NULL%1]'
